# Regenerate the "Demographics" sheet data for Gweru, Zimbabwe.
# (Loan given data source was re-pulled for the Gweru city dataset.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Helper: write a value as TEXT (shared string) even when it looks numeric,
# without leaving a quotePrefix flag or a new "Text" number-format style
# behind on the cell (mirrors how the source data was pasted in as text).
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# --- Column A: Respondent ID (now stored as text IDs from the new pull) ---
Set-TextValue $ws.Cells.Item(2, 1) "129938"
Set-TextValue $ws.Cells.Item(3, 1) "135128"
Set-TextValue $ws.Cells.Item(4, 1) "135430"
Set-TextValue $ws.Cells.Item(5, 1) "135525"
Set-TextValue $ws.Cells.Item(6, 1) "161237"
Set-TextValue $ws.Cells.Item(7, 1) "168729"

# --- Column I header + data: "Country of Residence" -> "City of Residence" ---
$ws.Cells.Item(1, 9).Value = "City of Residence"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "Gweru"
    $ws.Cells.Item($r, 10).Value = "Zimbabwe"
}

# --- Column B: Name (first name only now) ---
$ws.Cells.Item(2, 2).Value = "Gloria"
$ws.Cells.Item(3, 2).Value = "Patience"
$ws.Cells.Item(4, 2).Value = "Degreat"
$ws.Cells.Item(5, 2).Value = "Nnanna"
$ws.Cells.Item(6, 2).Value = "Igbire"
$ws.Cells.Item(7, 2).Value = "Iyamu"

# --- Column C: Email (new domain) ---
$ws.Cells.Item(2, 3).Value = "o1@xmail.com"
$ws.Cells.Item(3, 3).Value = "pa08@xmail.com"
$ws.Cells.Item(4, 3).Value = "e2@xmail.com"
$ws.Cells.Item(5, 3).Value = "p5@xmail.com"
$ws.Cells.Item(6, 3).Value = "be@xmail.com"
$ws.Cells.Item(7, 3).Value = "pas@xmail.com"

# Cursor / selection left on E13 (as in the authored workbook).
$ws.Range("E13").Select() | Out-Null
